# Appends " (Changed main)" after the first paragraph's existing text,
# as three separate runs: " (", "Changed main", ")" - matching the
# target OOXML diff exactly (the existing run is left untouched).

$d = $word.ActiveDocument

# Locate the paragraph that holds the sentence we need to extend.
$target = $d.Content
$found = $target.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$para = $d.Paragraphs(1).Range

# Pull the paragraph's own opening <w:p ...> tag (paraId/textId/rsid
# attributes) straight from the live document so the rebuilt paragraph
# keeps them unchanged.
$fullXml = $para.WordOpenXML
$nsDecl = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$pOpenTag = "<w:p" + $nsDecl + ">"
if ($fullXml -match "<w:p ([^>]*)>") {
    $pOpenTag = "<w:p " + $matches[1] + $nsDecl + ">"
} elseif ($fullXml -match "<w:p>") {
    $pOpenTag = "<w:p" + $nsDecl + ">"
}

$existingText = $para.Text
# Paragraph.Range.Text includes the trailing paragraph mark; strip it.
$existingText = $existingText.TrimEnd([char]13, [char]7)

function XmlEscape([string]$s) {
    $s = $s.Replace("&", "&amp;")
    $s = $s.Replace("<", "&lt;")
    $s = $s.Replace(">", "&gt;")
    return $s
}

$run1 = "<w:r><w:t>" + (XmlEscape $existingText) + "</w:t></w:r>"
$run2 = "<w:r><w:t xml:space=`"preserve`"> (</w:t></w:r>"
$run3 = "<w:r><w:t>Changed main</w:t></w:r>"
$run4 = "<w:r><w:t>)</w:t></w:r>"

$paragraphXml = $pOpenTag + $run1 + $run2 + $run3 + $run4 + "</w:p>"

$para.InsertXML($paragraphXml)
